# Fruta / hortaliza, semanal
# Adds a new week of Cereza price data (date 44939) for "Terminal La
# Palmera de La Serena" ahead of the existing "Lapins" rows, which are
# pushed down unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert six new blank rows above the current row 464 (the former
# "Lapins" rows 464-465 shift down to 470-471, untouched).
$ws.Range("464:469").Insert()

# Common columns shared by every new row.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$fecha     = 44939
$codreg    = 4
$tipo      = "Fruta"
$prodId    = 100103
$producto  = "Frutos de hueso (carozo)"
$catId     = 100103001
$categoria = "Cereza"

function Set-CerezaRow {
    param(
        [int]$row,
        [string]$variedad,
        [string]$calidad,
        [double]$volumen,
        [double]$precioMin,
        [double]$precioMax,
        [double]$precioProm,
        [string]$unidad,
        [string]$origen,
        [double]$precioKg,
        [double]$kgUnidad
    )

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $prodId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $catId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-CerezaRow 464 "Bing"        "Especial" 340 12000 13000 12500 "$/bandeja 10 kilos" "Provincia de Curicó" 1250 10
Set-CerezaRow 465 "Bing"        "Primera"  300 10000 11000 10500 "$/bandeja 10 kilos" "Provincia de Curicó" 1050 10
Set-CerezaRow 466 "Bing"        "Segunda"  400  8000  9000  8500 "$/bandeja 10 kilos" "Provincia de Curicó"  850 10
Set-CerezaRow 467 "Sweet Heart" "Especial" 300 15000 16000 15500 "$/caja 15 kilos"    "Provincia de Curicó" 1033 15
Set-CerezaRow 468 "Sweet Heart" "Primera"  400 13000 14000 13500 "$/caja 15 kilos"    "Provincia de Curicó"  900 15
Set-CerezaRow 469 "Sweet Heart" "Segunda"  300 11000 12000 11500 "$/caja 15 kilos"    "Provincia de Curicó"  767 15
